$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exp 4 (row 6) result update
$ws.Range("N6").Value = 266.82029999999997
$ws.Range("O6").Value = 270

# New Exp 5 result (row 7)
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "DDPG"
$ws.Range("C7").Value = 2500
$ws.Range("D7").Value = "0.001 / 0.001"
$ws.Range("E7").Value = 64
$ws.Range("F7").Value = 500000
$ws.Range("G7").Value = "Adam"
$ws.Range("H7").Value = "nan"
$ws.Range("I7").Value = "nan"
$ws.Range("J7").Value = 0.99
$ws.Range("K7").Value = 4
$ws.Range("L7").Value = 4
$ws.Range("M7").Value = 0.001
$ws.Range("O7").Value = 300

$ws.Rows.Item(7).RowHeight = 21.75

$ws.Range("N7").Select() | Out-Null
